$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 16.3.1: add a 2023 data column (E) alongside the existing 2018 column (D),
# bold the "Urbanisation / Местность / Жерлери" section header row, add the
# missing urban/rural 2023 breakdown ("-") and update the footnote to mention
# the 2023 survey as well as 2018.
# ---------------------------------------------------------------------------

# --- Column E header year (2023), formatted like D4 (2018) -----------------
$ws.Range("D4").Copy() | Out-Null
$ws.Range("E4").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("E4").Value = 2023

# --- Column E overall value for the indicator row, formatted like D5 -------
$ws.Range("D5").Copy() | Out-Null
$ws.Range("E5").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("E5").Value = 38

# --- Row 6 ("Жерлери" / "Местность" / "Urbanisation") becomes bold ---------
# Give the new E6 cell the same base formatting as D6 first, then bold the
# whole row so the new cell ends up on the same font as the rest of the row.
$ws.Range("D6").Copy() | Out-Null
$ws.Range("E6").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("A6:E6").Font.Bold = $true

# --- Column E urban/rural rows: no 2023 breakdown available ("-") ----------
$ws.Range("D7").Copy() | Out-Null
$ws.Range("E7").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("E7").HorizontalAlignment = -4152      # xlRight
$ws.Range("E7").Value = "-"

$ws.Range("D8").Copy() | Out-Null
$ws.Range("E8").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("E8").HorizontalAlignment = -4152      # xlRight
$ws.Range("E8").Value = "-"

# --- Footnote row: mention the 2023 survey alongside 2018 -------------------
$ws.Range("A9").Value = " Көп көрсөткүчтүү кластердик изилдөөнүн маалыматтары боюнча, 2018-ж., 2023-ж."
$ws.Range("B9").Value = "По данным кластерного обследования по многим показателям, 2018г., 2023г."
$ws.Range("C9").Value = "According to the cluster survey in many respects, 2018, 2023."
